$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price / volume(1h) figures.
# These columns hold numeric-looking text (e.g. "5.530", "0.49%") where the
# exact textual representation (trailing zeros, % suffix) must be preserved,
# so force each cell to Text format before writing the new value - otherwise
# Excel would silently convert it to a floating point number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.49%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.45%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.530"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.81%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08146"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.65%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.27%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9758"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.41%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.1115"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-3.98%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1891"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.05%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "10.16"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-19.37%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1003"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.91%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04787"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1058"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.91%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001254"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.46%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04106"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.02%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005975"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.58%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.343"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.92%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.429"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.26%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.70%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.54%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.57%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.81%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004383"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.40%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001277"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.23%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003731"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.24%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02690"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.05%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05655"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007608"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.52%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.79%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007529"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.89%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001954"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.14%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008291"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-6.89%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007017"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.33%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000748"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.33%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005786"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.44%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002514"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.65%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003523"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-19.06%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.33%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001995"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.33%"

# Coin / Link columns shifted as the ranking list was refreshed
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B8").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C8").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B10").Value = "MCDex"
$ws.Range("C10").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
